# Weekly update: insert a new price record for "Arveja Verde" at row 145,
# pushing the existing rows 145-188 down to 146-189.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 145 (shifts rows 145:188 down to 146:189).
$ws.Rows("145:145").Insert()

# Populate the newly inserted row 145 with the new record.
$ws.Range("A145").Value = 6
$ws.Range("B145").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C145").Value = "Metropolitana"
$ws.Range("D145").Value = 44551
$ws.Range("E145").Value = 13
$ws.Range("F145").Value = 100112022
$ws.Range("G145").Value = "Arveja Verde"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 300
$ws.Range("K145").Value = 10000
$ws.Range("L145").Value = 12000
$ws.Range("M145").Value = 10800
$ws.Range("N145").Value = "`$/saco 25 kilos"
$ws.Range("O145").Value = "Carahue"
$ws.Range("P145").Value = 432
$ws.Range("Q145").Value = 25
$ws.Range("R145").Value = "Hortaliza"
